# Generate Report for Handoff
# Refresh the "Ready for handoff" rows (4-7) across the Overview, zh-cn and
# de-de sheets: the Priority moves from "low" to "ht" and the handoff
# timestamps are bumped to reflect the freshly generated handoff package.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-08-15 22:29:42"

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-08-15 22:29:47"

    $overview.Range("G$row").Value = "2016-08-15 22:29:47"
}
